$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$newValues = @(
    "Comportamentos de interação social malsucedidos ",
    "Desconforto em situações sociais ",
    "Incapacidade de comunicar uma sensação satisfatória de envolvimento social (p. ex., pertinência, cuidado, interesse, história compartilhada) ",
    "Incapacidade de receber uma sensação satisfatória de envolvimento social (p. ex., pertinência, cuidado, interesse, história compartilhada) ",
    "Interação disfuncional com outras pessoas ",
    "Relato familiar de mudança na interação (p. ex., estilo, padrão) "
)

$startRow = 332
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "###"
    $ws.Cells.Item($r, 2).Value = 115
    $ws.Cells.Item($r, 3).Value = "%%%"
    $ws.Cells.Item($r, 4).Value = $newValues[$i]
    $ws.Cells.Item($r, 5).Value = "$$$"
}

$ws.Range("A331:A337").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 313
